$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.770.23"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "1.700.28"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'316.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").Value = "'0.3929"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").Value = "'0.4039"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'1.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").Value = "'53.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "

$ws.Range("D11").Value = "'1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "'0.08891"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").Value = "'7.248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").Value = "'23.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").Value = "'8.033"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.30%  "

$ws.Range("D16").Value = "'0.00001326"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").Value = "1.703.46"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "'100.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "'0.07016"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'19.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "'7.050"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'14.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.57%  "

$ws.Range("D24").Value = "24.755.08"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").Value = "'3.247"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.08%  "

$ws.Range("D26").Value = "'2.355"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").Value = "'161.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("D29").Value = "'136.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("D31").Value = "'7.762"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").Value = "'0.08742"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").Value = "'1.074"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.57%  "

$ws.Range("D34").Value = "'7.217"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.61%  "

$ws.Range("D35").Value = "'11.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").Value = "'1.961"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").Value = "'0.2749"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("E38").Value = "  -2.87%  "

$ws.Range("D39").Value = "'0.09192"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "'0.02736"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").Value = "'1.464"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").Value = "'0.7684"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'15.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("D44").Value = "'0.7178"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "

$ws.Range("D45").Value = "'2.574"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("D46").Value = "'4.220"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").Value = "'140.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("D49").Value = "'1.309"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").Value = "'0.07982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").Value = "'90.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.52%  "
